$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Left-hand table (REQUISITO 23: Registro de usuario / Inicio de sesion)
# ---------------------------------------------------------------

# Row 4: task name text changed, D5 value later adjusted
$ws.Range("B4").Value2 = "Video Tutoriales/formación Python"
$ws.Range("F4").Value2 = "Video Tutoriales"

# Row 5
$ws.Range("D5").Value2 = 24

# Row 6
$ws.Range("D6").Value2 = 20

# Row 7
$ws.Range("D7").Value2 = 20

# Row 8
$ws.Range("D8").Value2 = 35
$ws.Range("F8").Value2 = "Código"

# Row 9: D9 and H9 become blank (keep border style), E9 becomes a plain number
$ws.Range("D9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("E9").Value2 = 30

# Recompute totals for the first table
$ws.Range("D11").Formula = "=D4+D5+D6+D8+D7+D9+D10"
$ws.Range("H11").Formula = "=SUM(H4:H10)"
$ws.Range("E11").Formula = "=D11+H11+E9"

# ---------------------------------------------------------------
# Second table (REQUISITO 25: Cambiar Contraseña)
# ---------------------------------------------------------------

# Row 15 keeps "Video Tutoriales" / "Miguel" / 20 but used to reference the
# shared "Miguel y Juan Pablo" text - now split per person
$ws.Range("C15").Value2 = "Miguel"

# Row 16
$ws.Range("C16").Value2 = "Juan Pablo"

# Prepare new rows 17-20 by copying the format from an existing data row,
# then row 21 (new TOTAL row) copies the format from the old TOTAL row (17)
$ws.Range("B17:D17").Copy() | Out-Null
$ws.Range("B21:D21").PasteSpecial(-4122) | Out-Null

$ws.Range("B15:D15").Copy() | Out-Null
$ws.Range("B17:D20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new content rows
$ws.Range("B17").Value2 = "Integración "
$ws.Range("C17").Value2 = "Juan Pablo"
$ws.Range("D17").Value2 = 30

$ws.Range("B18").Value2 = "Desarrollo Interfaz"
$ws.Range("C18").Value2 = "Juan Pablo"
$ws.Range("D18").Value2 = 13

$ws.Range("B19").Value2 = "Desarrollo Código Introducir Nueva Contraseña"
$ws.Range("C19").Value2 = "Miguel"
$ws.Range("D19").Value2 = 45

$ws.Range("B20").Value2 = "Desarrollo Código Requisitos Cambio Contraseña"
$ws.Range("C20").Value2 = "Juan Pablo"
$ws.Range("D20").Value2 = 55

# New TOTAL row
$ws.Range("B21").Value2 = "TOTAL"
$ws.Range("D21").Formula = "=D15+D16+D19+D17+D20+D18"

# ---------------------------------------------------------------
# Cosmetic sheet-level tweaks
# ---------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 44
$ws.Range("F20").Select() | Out-Null
